$wb = $excel.ActiveWorkbook

$wsTestList = $wb.Worksheets.Item("Test_List")
$wsPortland = $wb.Worksheets.Item("Portland")
$wsMiami    = $wb.Worksheets.Item("Miami")

# Template cell already carrying the "fontId=2 / General" cell style (same style as
# Portland!A2 originally). C2 is never touched by this script, so it stays a stable
# paste-special (formats only) source for the whole run -- restyled Test_ID cells land
# on the same font as the rest of the data column instead of inheriting a new font record.
$fmtSource = $wsPortland.Range("C2")

# ---- Test_List: Test_ID 1.0 -> 1.1 (reuses existing General/fontId=2 style) ----
$fmtSource.Copy()
$wsTestList.Range("A2").PasteSpecial(-4122)
$wsTestList.Range("A2").Value = 1.1

# ---- Portland: every Test_ID column gets an explicit numeric format ----
# Sub-id rows (value 1.0) -> one decimal place "0.0"; becomes new style index 3.
$fmtSource.Copy()
$wsPortland.Range("A2").PasteSpecial(-4122)
$wsPortland.Range("A2").NumberFormat = "0.0"
$wsPortland.Range("A2").Value = 1.1

# Integer rows (values 2.0/3.0/4.0) -> whole-number "0"; becomes new style index 4.
$fmtSource.Copy()
$wsPortland.Range("A3:A5").PasteSpecial(-4122)
$wsPortland.Range("A3:A5").NumberFormat = "0"

$fmtSource.Copy()
$wsPortland.Range("A6").PasteSpecial(-4122)
$wsPortland.Range("A6").NumberFormat = "0.0"
$wsPortland.Range("A6").Value = 1.1

$fmtSource.Copy()
$wsPortland.Range("A7:A9").PasteSpecial(-4122)
$wsPortland.Range("A7:A9").NumberFormat = "0"

$fmtSource.Copy()
$wsPortland.Range("A10").PasteSpecial(-4122)
$wsPortland.Range("A10").NumberFormat = "0.0"
$wsPortland.Range("A10").Value = 1.1

$fmtSource.Copy()
$wsPortland.Range("A11:A13").PasteSpecial(-4122)
$wsPortland.Range("A11:A13").NumberFormat = "0"

# ---- Miami: only the sub-id rows (value 1.0) pick up the "0.0" format ----
$fmtSource.Copy()
$wsMiami.Range("A2").PasteSpecial(-4122)
$wsMiami.Range("A2").NumberFormat = "0.0"
$wsMiami.Range("A2").Value = 1.1

$fmtSource.Copy()
$wsMiami.Range("A5").PasteSpecial(-4122)
$wsMiami.Range("A5").NumberFormat = "0.0"
$wsMiami.Range("A5").Value = 1.1

$fmtSource.Copy()
$wsMiami.Range("A10").PasteSpecial(-4122)
$wsMiami.Range("A10").NumberFormat = "0.0"
$wsMiami.Range("A10").Value = 1.1

$excel.CutCopyMode = 0
